$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns (antecedent_len, consequent_len, combo_len) ---
$ws.Range("H1").Value = "antecedent_len"
$ws.Range("I1").Value = "consequent_len"
$ws.Range("J1").Value = "combo_len"

# Match the header formatting (bold font, centered, bordered) used by A1:G1
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2 ---
$ws.Range("A2").Value = "frozenset({'ACCESS_NETWORK_STATE'})"
$ws.Range("B2").Value = "frozenset({'INTERNET'})"
$ws.Range("C2").Value = 0.9102564102564102
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1.012987012987013
$ws.Range("F2").Value = 0.01166995397764625
$ws.Range("G2").Value = "inf"
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 2

# --- Row 3 ---
$ws.Range("A3").Value = "frozenset({'ACCESS_NETWORK_STATE', 'WAKE_LOCK'})"
$ws.Range("B3").Value = "frozenset({'INTERNET'})"
$ws.Range("C3").Value = 0.8333333333333334
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1.012987012987013
$ws.Range("F3").Value = 0.01068376068376065
$ws.Range("G3").Value = "inf"
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 3

# --- Row 4 ---
$ws.Range("A4").Value = "frozenset({'WAKE_LOCK', 'INTERNET'})"
$ws.Range("B4").Value = "frozenset({'ACCESS_NETWORK_STATE'})"
$ws.Range("C4").Value = 0.8333333333333334
$ws.Range("D4").Value = 0.9848484848484848
$ws.Range("E4").Value = 1.081946222791293
$ws.Range("F4").Value = 0.0631163708086786
$ws.Range("G4").Value = 5.92307692307693
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 3

# --- Row 5 ---
$ws.Range("A5").Value = "frozenset({'ACCESS_NETWORK_STATE'})"
$ws.Range("B5").Value = "frozenset({'WAKE_LOCK'})"
$ws.Range("C5").Value = 0.8333333333333334
$ws.Range("D5").Value = 0.9154929577464788
$ws.Range("E5").Value = 1.081946222791293
$ws.Range("F5").Value = 0.0631163708086786
$ws.Range("G5").Value = 1.820512820512821
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 2

# --- Row 6 (new) ---
$ws.Range("A6").Value = "frozenset({'ACCESS_NETWORK_STATE'})"
$ws.Range("B6").Value = "frozenset({'WAKE_LOCK', 'INTERNET'})"
$ws.Range("C6").Value = 0.8333333333333334
$ws.Range("D6").Value = 0.9154929577464788
$ws.Range("E6").Value = 1.081946222791293
$ws.Range("F6").Value = 0.0631163708086786
$ws.Range("G6").Value = 1.820512820512821
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 2
$ws.Range("J6").Value = 3

# --- Row 7 (new) ---
$ws.Range("A7").Value = "frozenset({'INTERNET'})"
$ws.Range("B7").Value = "frozenset({'WAKE_LOCK'})"
$ws.Range("C7").Value = 0.8461538461538461
$ws.Range("D7").Value = 0.8571428571428571
$ws.Range("E7").Value = 1.012987012987013
$ws.Range("F7").Value = 0.01084812623274156
$ws.Range("G7").Value = 1.076923076923077
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 2

Write-Output "edit applied"
